$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-6: update financial figures with corrected values
$cellValues = @{
    "D2" = 4569
    "E2" = 47
    "F2" = 47
    "G2" = 1402
    "H2" = 1132
    "I2" = 1391
    "J2" = -260
    "K2" = 16383
    "L2" = 12686
    "M2" = 3697
    "N2" = 3198
    "O2" = 499
    "P2" = 287
    "Q2" = 664
    "R2" = 1881
    "S2" = -1946
    "T2" = 66
    "U2" = 598
    "V2" = 2980
    "W2" = 1.03
    "X2" = 24.77
    "Y2" = 52.46
    "Z2" = 5.73
    "AA2" = 343.17
    "AB2" = 871.54
    "AC2" = 2633
    "AD2" = 0.73
    "AE2" = 6074
    "AF2" = 0.32
    "AG2" = 0
    "AH2" = 0
    "AI2" = 0
    "AJ2" = 52866795
    "D3" = 4373
    "E3" = 187
    "F3" = 187
    "G3" = 327
    "H3" = 204
    "I3" = 308
    "J3" = -104
    "K3" = 12037
    "L3" = 8462
    "M3" = 3575
    "N3" = 3526
    "O3" = 48
    "P3" = 323
    "Q3" = 830
    "R3" = -121
    "S3" = -1830
    "T3" = 29
    "U3" = 801
    "V3" = 675
    "W3" = 4.28
    "X3" = 4.67
    "Y3" = 9.16
    "Z3" = 1.44
    "AA3" = 236.73
    "AB3" = 905.76
    "AC3" = 579
    "AD3" = 5.7
    "AE3" = 5891
    "AF3" = 0.5600000000000001
    "AG3" = 0
    "AH3" = 0
    "AI3" = 0
    "AJ3" = 60070918
    "D4" = 2430
    "E4" = -118
    "F4" = -118
    "G4" = 56
    "H4" = -332
    "I4" = -251
    "J4" = -81
    "K4" = 8074
    "L4" = 4819
    "M4" = 3255
    "N4" = 3333
    "O4" = -78
    "P4" = 386
    "Q4" = 163
    "R4" = 1589
    "S4" = -1146
    "T4" = 10
    "U4" = 153
    "V4" = 520
    "W4" = -4.85
    "X4" = -13.65
    "Y4" = -7.32
    "Z4" = -3.3
    "AA4" = 148.06
    "AB4" = 733.74
    "AC4" = -412
    "AD4" = -6.86
    "AE4" = 4600
    "AF4" = 0.61
    "AG4" = 0
    "AH4" = 0
    "AI4" = 0
    "AJ4" = 72677207
    "D5" = 2828
    "E5" = 104
    "F5" = 104
    "G5" = 155
    "H5" = 69
    "I5" = 86
    "J5" = -16
    "K5" = 8862
    "L5" = 5525
    "M5" = 3337
    "N5" = 3429
    "O5" = -32
    "P5" = 386
    "Q5" = 148
    "R5" = -1302
    "S5" = 507
    "T5" = 49
    "U5" = 99
    "V5" = 1098
    "W5" = 3.67
    "X5" = 2.46
    "Y5" = 2.54
    "Z5" = 0.82
    "AA5" = 165.57
    "AB5" = 759.71
    "AC5" = 118
    "AD5" = 18.96
    "AE5" = 4733
    "AF5" = 0.47
    "AG5" = 0
    "AH5" = 0
    "AI5" = 0
    "AJ5" = 72677207
    "D6" = 2904
    "E6" = -842
    "F6" = -842
    "G6" = -1011
    "H6" = -1033
    "I6" = -1017
    "K6" = 7599
    "L6" = 5342
    "M6" = 2256
    "N6" = 2369
    "P6" = 393
    "Q6" = -567
    "R6" = 834
    "S6" = -451
    "T6" = 195
    "U6" = -762
    "V6" = 926
    "W6" = -29
    "X6" = -35.56
    "Y6" = -35.07
    "Z6" = -12.55
    "AA6" = 236.75
    "AB6" = 477.13
    "AC6" = -1397
    "AD6" = -1.51
    "AE6" = 3203
    "AF6" = 0.66
    "AG6" = 0
    "AH6" = 0
    "AI6" = 0
    "AJ6" = 74166668
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}

# Rows 7-9: clear erroneous data cells (D:AJ), keep only A-C identifiers
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
